$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": PIEDRA SINTERIZADA (col L) sale for client DDH S.A.S. (row 11) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L11").Value = 253.44
$ws1.Range("L24").Value = "1 de 22"

# --- Sheet "VENTA MENSUAL": agosto (col F) sale for client DDH S.A.S. (row 11) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F11").Value = 253.44
$ws2.Range("F24").Value = 3757.45

# --- Sheet "CUMPLIMIENTO MENSUAL": PIEDRA SINTERIZADA group (row 15) and TOTAL (row 19) ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D15").Value = 253.44
$ws3.Range("E15").Value = 2247.57
$ws3.Range("F15").Value = 0.1013350606355033

$ws3.Range("D19").Value = 3757.45
$ws3.Range("E19").Value = 51265.71386304604
$ws3.Range("F19").Value = 0.06828851225917111
